$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.8109359125033627
$ws.Range("D2").Value = 0.4230429764014652

$ws.Range("C3").Value = 1.347037848281134
$ws.Range("D3").Value = 0.1868781134112933

$ws.Range("C4").Value = 1.268430223226455
$ws.Range("D4").Value = 0.2132608143202759

$ws.Range("C5").Value = 0.02890456400193659
$ws.Range("D5").Value = 0.977109706868786

$ws.Range("C6").Value = 0.6371744905960612
$ws.Range("D6").Value = 0.5282786689916827

$ws.Range("C7").Value = 0.7859003828187799
$ws.Range("D7").Value = 0.4373689178495577

$ws.Range("C8").Value = -0.6496292607672908
$ws.Range("D8").Value = 0.520296798542649

$ws.Range("C9").Value = -0.003706835317814666
$ws.Range("D9").Value = 0.9970640446400343

$ws.Range("C10").Value = -1.52837306847804
$ws.Range("D10").Value = 0.1356709250159636

$ws.Range("C11").Value = -1.420207278379069
$ws.Range("D11").Value = 0.1646528726724199
